$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: simple value replacements in the first four rows ---
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
$t.Cell(4, 1).Range.Text = "211"

# --- Step 2: insert three new rows right after row 4 (i.e. before row 5) ---
# Rows.Add(beforeRow) inserts immediately above beforeRow, so add them in
# reverse order to end up with 0.00002, 0.00060, 0.00018 (top to bottom).
$anchorRow = $t.Rows.Item(5)

$newRow3 = $t.Rows.Add($anchorRow)
$newRow3.Cells.Item(1).Range.Text = "0.00018"

$newRow2 = $t.Rows.Add($anchorRow)
$newRow2.Cells.Item(1).Range.Text = "0.00060"

$newRow1 = $t.Rows.Add($anchorRow)
$newRow1.Cells.Item(1).Range.Text = "0.00002"

# --- Step 3: update the rows that followed the original row 4 ---
# After the 3 insertions, the former rows 6,7,8,9 are now at 9,10,11,12.
$t.Cell(9, 1).Range.Text = "0.00034"
$t.Cell(10, 1).Range.Text = "0.00040"
$t.Cell(11, 1).Range.Text = "0.00043"
$t.Cell(12, 1).Range.Text = "0.04589"

# --- Step 4: delete the three rows that followed (former rows 10,11,12,
# now at positions 13,14,15). Delete from the bottom up so indices stay valid.
$t.Rows.Item(15).Delete()
$t.Rows.Item(14).Delete()
$t.Rows.Item(13).Delete()

# --- Step 5: collapse the three large tab-separated summary rows down to a
# single value each. These rows are unaffected by the earlier insert/delete
# (3 added, 3 removed before them), so they remain rows 44, 45, 46.
$t.Cell(44, 1).Range.Text = "99.92"
$t.Cell(45, 1).Range.Text = "0.05"
$t.Cell(46, 1).Range.Text = "56"
